# TODO.docx fix: shrink the bullet-list font size (32 half-points -> 24,
# i.e. 16pt -> 12pt) and give every run an explicit complex-script size
# (szCs 18 half-points = 9pt), matching the "fix not found concept" commit.
# The document's "_GoBack" bookmark (Word's "last edit position" marker)
# also needs to end up right after "graphTab" instead of after "Check w".

$d = $word.ActiveDocument

# --- 1. Move the "_GoBack" bookmark to right after "graphTab" -------------
# Remove it from its current location (end of the "Check w" paragraph).
$old = $d.Bookmarks("_GoBack")
$old.Delete()

# Locate "graphTab" and collapse the found range to its end point.
$target = $d.Content
$target.Find.Execute("graphTab", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$target.Collapse(0)

# A genuinely zero-length range can't be handed to Bookmarks.Add directly
# here, so stage a 1-char placeholder, bookmark across it, then delete the
# placeholder - the bookmark collapses back to a single point in place.
$target.InsertBefore("X")
$d.Bookmarks.Add("_GoBack", $target)
$target.Text = ""

# --- 2. Shrink the font everywhere except the trailing empty paragraph ----
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -lt $paraCount; $i++) {
    $font = $d.Paragraphs($i).Range.Font
    $font.Size = 12
    $font.SizeBi = 9
}
